$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion message text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.41 = 30151.7 pesos`n✅ 30151.7 pesos = 7.38 = 967.75 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 135
$ws2.Range("O10").Value = 4070.48
$ws2.Range("N12").Value = 4084
$ws2.Range("O12").Value = 131.08
